# Refactor synthetic array: insert a new "statut_name" column (C) between
# "statut_label" (B) and "NCTId" (old C, now D), shifting all subsequent
# columns one to the right. Populate the new column with a human-readable
# label derived from the existing "statut_label" (noir/vert) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; everything from the old C onward shifts right.
$ws.Columns(3).Insert()

# Header for the new column.
$ws.Range("C1").Value = "statut_name"

$statusMap = @{
    "noir" = "pas de résultat ni de publication"
    "vert" = "résultat et / ou publication posté dans les 12 mois"
}

$lastRow = $ws.Cells(1, 1).Worksheet.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $label = $ws.Cells($row, 2).Value2
    if ($statusMap.ContainsKey($label)) {
        $ws.Cells($row, 3).Value = $statusMap[$label]
    }
}
